$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-10 (row 23)
$ws.Range("B23").Value = 6305
$ws.Range("C23").Value = 997
$ws.Range("D23").Value = 5867038
$ws.Range("E23").Value = 930.5373513084853
$ws.Range("F23").Value = 8.184625943719981
$ws.Range("G23").Value = 3.746097814776284
$ws.Range("H23").Value = 25.71560686903742
